$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 477; this shifts the existing rows 477-502
# down to 478-503 (matching every downstream row-shift in the diff) and
# carries the date number-format from the row above into the new row.
$ws.Rows.Item(477).Insert()

# Populate the newly inserted row 477 with this week's data point.
$ws.Range("A477").Value = 11
$ws.Range("B477").Value = "Vega Monumental Concepción"
$ws.Range("C477").Value = "Bíobío"
$ws.Range("D477").Value = 45106
$ws.Range("E477").Value = 8
$ws.Range("F477").Value = "Fruta"
$ws.Range("G477").Value = 100102
$ws.Range("H477").Value = "Cítricos"
$ws.Range("I477").Value = 100102005
$ws.Range("J477").Value = "Naranja"
$ws.Range("K477").Value = "Lane Late"
$ws.Range("L477").Value = "Primera"
$ws.Range("M477").Value = 220
$ws.Range("N477").Value = 9000
$ws.Range("O477").Value = 9500
$ws.Range("P477").Value = 9227
$ws.Range("Q477").Value = "$/bandeja 15 kilos granel"
$ws.Range("R477").Value = "Región de O'Higgins"
$ws.Range("S477").Value = 615
$ws.Range("T477").Value = 15
